$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff swaps the data between row 2 and row 3 for columns D (Fecha),
# J (Volumen), K (Precio mínimo), M (Precio promedio ponderado) and P (Precio $/Kg).
# Column L (Precio máximo) stays the same (10000) for both rows.

$ws.Range("D2").Value = 44804
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 9500
$ws.Range("M2").Value = 9750
$ws.Range("P2").Value = 542

$ws.Range("D3").Value = 44714
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 9000
$ws.Range("M3").Value = 9500
$ws.Range("P3").Value = 528
